$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37729
$ws.Range("D2").Value = 54565565
$ws.Range("C3").Value = 90956
$ws.Range("D3").Value = 133332235
$ws.Range("C4").Value = 31175
$ws.Range("D4").Value = 46168849
$ws.Range("C5").Value = 8691
$ws.Range("D5").Value = 12917063
$ws.Range("C6").Value = 1993
$ws.Range("D6").Value = 2962006
$ws.Range("C12").Value = 41315
$ws.Range("D12").Value = 56056715
$ws.Range("C13").Value = 9650
$ws.Range("D13").Value = 13956958
$ws.Range("C14").Value = 25939
$ws.Range("D14").Value = 38042786
$ws.Range("C15").Value = 8309
$ws.Range("D15").Value = 12331324
$ws.Range("C20").Value = 10219
$ws.Range("D20").Value = 13530625
$ws.Range("C21").Value = 13375
$ws.Range("D21").Value = 19312792
$ws.Range("C22").Value = 31650
$ws.Range("D22").Value = 46445820
$ws.Range("C23").Value = 10216
$ws.Range("D23").Value = 15186678
$ws.Range("C24").Value = 2638
$ws.Range("D24").Value = 3922182
$ws.Range("C27").Value = 11684
$ws.Range("D27").Value = 15607194
$ws.Range("C28").Value = 7640
$ws.Range("D28").Value = 11066617
$ws.Range("C29").Value = 22474
$ws.Range("D29").Value = 32989049
$ws.Range("C30").Value = 7813
$ws.Range("D30").Value = 11627133
$ws.Range("C31").Value = 1958
$ws.Range("D31").Value = 2921499
$ws.Range("C34").Value = 8305
$ws.Range("D34").Value = 10970625
$ws.Range("C35").Value = 3245
$ws.Range("D35").Value = 4684194
$ws.Range("C36").Value = 7823
$ws.Range("D36").Value = 11424434
$ws.Range("C37").Value = 3177
$ws.Range("D37").Value = 4708461
$ws.Range("C42").Value = 17229
$ws.Range("D42").Value = 24910878
$ws.Range("C43").Value = 51082
$ws.Range("D43").Value = 74886786
$ws.Range("C44").Value = 19007
$ws.Range("D44").Value = 28232943
$ws.Range("C50").Value = 16695
$ws.Range("D50").Value = 22226805
$ws.Range("C51").Value = 2020
$ws.Range("D51").Value = 2929471
$ws.Range("C52").Value = 6893
$ws.Range("D52").Value = 10132824
$ws.Range("C53").Value = 2348
$ws.Range("D53").Value = 3506918
$ws.Range("C54").Value = 755
$ws.Range("D54").Value = 1127805
$ws.Range("C55").Value = 186
$ws.Range("D55").Value = 275833
$ws.Range("C57").Value = 6972
$ws.Range("D57").Value = 9586294
$ws.Range("C58").Value = 945
$ws.Range("D58").Value = 1387079
$ws.Range("C59").Value = 2374
$ws.Range("D59").Value = 3519837
$ws.Range("C60").Value = 943
$ws.Range("D60").Value = 1404001
$ws.Range("C61").Value = 322
$ws.Range("D61").Value = 482758
$ws.Range("C62").Value = 104
$ws.Range("D62").Value = 155850
$ws.Range("C64").Value = 1392
$ws.Range("D64").Value = 1959206
$ws.Range("C65").Value = 15351
$ws.Range("D65").Value = 22175263
$ws.Range("C66").Value = 44662
$ws.Range("D66").Value = 65357707
$ws.Range("C67").Value = 15695
$ws.Range("D67").Value = 23324686
$ws.Range("C68").Value = 4567
$ws.Range("D68").Value = 6802792
$ws.Range("C73").Value = 15076
$ws.Range("D73").Value = 19877419
$ws.Range("C74").Value = 51350
$ws.Range("D74").Value = 74726300
$ws.Range("C75").Value = 145944
$ws.Range("D75").Value = 215009719
$ws.Range("C76").Value = 63599
$ws.Range("D76").Value = 94771307
$ws.Range("C77").Value = 20330
$ws.Range("D77").Value = 30375331
$ws.Range("C78").Value = 4814
$ws.Range("D78").Value = 7190043
$ws.Range("C79").Value = 264
$ws.Range("D79").Value = 391170
$ws.Range("C85").Value = 50778
$ws.Range("D85").Value = 69076673
$ws.Range("C86").Value = 4596
$ws.Range("D86").Value = 6658911
$ws.Range("C87").Value = 11556
$ws.Range("D87").Value = 16976869
$ws.Range("C89").Value = 1343
$ws.Range("D89").Value = 2006989
$ws.Range("C93").Value = 5409
$ws.Range("D93").Value = 7271613
$ws.Range("C94").Value = 1595
$ws.Range("D94").Value = 2297432
$ws.Range("C95").Value = 5161
$ws.Range("D95").Value = 7600243
$ws.Range("C96").Value = 1939
$ws.Range("D96").Value = 2888437
$ws.Range("C101").Value = 3557
$ws.Range("D101").Value = 4707764
$ws.Range("C102").Value = 601
$ws.Range("D102").Value = 895164
$ws.Range("C107").Value = 10748
$ws.Range("D107").Value = 15591962
$ws.Range("C108").Value = 29189
$ws.Range("D108").Value = 42886096
$ws.Range("C109").Value = 9773
$ws.Range("D109").Value = 14533150
$ws.Range("C110").Value = 2682
$ws.Range("D110").Value = 3999207
$ws.Range("C114").Value = 9788
$ws.Range("D114").Value = 12929894
$ws.Range("C115").Value = 30431
$ws.Range("D115").Value = 43883084
$ws.Range("C116").Value = 66112
$ws.Range("D116").Value = 96753676
$ws.Range("C117").Value = 21356
$ws.Range("D117").Value = 31738440
$ws.Range("C118").Value = 6063
$ws.Range("D118").Value = 9032521
$ws.Range("C119").Value = 1120
$ws.Range("D119").Value = 1673771
$ws.Range("C124").Value = 25829
$ws.Range("D124").Value = 34499068
$ws.Range("C125").Value = 35953
$ws.Range("D125").Value = 51888527
$ws.Range("C126").Value = 76738
$ws.Range("D126").Value = 112212982
$ws.Range("C127").Value = 23830
$ws.Range("D127").Value = 35366409
$ws.Range("C128").Value = 6390
$ws.Range("D128").Value = 9495738
$ws.Range("C129").Value = 1234
$ws.Range("D129").Value = 1835411
$ws.Range("C133").Value = 31787
$ws.Range("D133").Value = 42207889
$ws.Range("C134").Value = 13204
$ws.Range("D134").Value = 19112130
$ws.Range("C135").Value = 32309
$ws.Range("D135").Value = 47454152
$ws.Range("C136").Value = 11468
$ws.Range("D136").Value = 17040042
$ws.Range("C137").Value = 2954
$ws.Range("D137").Value = 4403714
$ws.Range("C138").Value = 500
$ws.Range("D138").Value = 743990
$ws.Range("C141").Value = 10802
$ws.Range("D141").Value = 14405556
$ws.Range("C142").Value = 35012
$ws.Range("D142").Value = 50557490
$ws.Range("C143").Value = 81159
$ws.Range("D143").Value = 118908615
$ws.Range("C144").Value = 24318
$ws.Range("D144").Value = 36131155
$ws.Range("C145").Value = 6383
$ws.Range("D145").Value = 9524067
$ws.Range("C146").Value = 1433
$ws.Range("D146").Value = 2131730
$ws.Range("C149").Value = 29174
$ws.Range("D149").Value = 39358010
